$d = $word.ActiveDocument

# --- 1. Letter date: "July 8, 2021" -> "October 9th, 2021" (with "th" superscript) ---
$d.Content.Find.Execute("July 8", $true, $false, $false, $false, $false, $true, 1, $false, "October 9", 2) | Out-Null

$dateRange = $d.Range(0, 0)
$dateRange.Find.Execute("October 9", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dateRange.Collapse(0)  # wdCollapseEnd
$dateRange.InsertAfter("th")
$dateRange.Font.Superscript = $true

# --- 2. Recipient name block ---
$d.Content.Find.Execute("Antonia Hamilton", $true, $false, $false, $false, $false, $true, 1, $false, "Anastasia Efklides", 2) | Out-Null
$d.Content.Find.Execute("Department of Psychology", $true, $false, $false, $false, $false, $true, 1, $false, "School of Psychology", 2) | Out-Null
$d.Content.Find.Execute("University College London", $true, $false, $false, $false, $false, $true, 1, $false, "Aristotle University of Thessaloniki", 2) | Out-Null
$d.Content.Find.Execute("London, UK", $true, $false, $false, $false, $false, $true, 1, $false, "Thessaloniki, Greece", 2) | Out-Null

# --- 3. Salutation: "Dear Dr. Hamilton" -> "Dear Dr. Efklides" ---
$d.Content.Find.Execute("Hamilton", $true, $false, $false, $false, $false, $true, 1, $false, "Efklides", 2) | Out-Null

# --- 4. Journal name: "the Quarterly Journal of Experimental Psychology" -> "Metacognition and Learning" ---
# (two occurrences; the first one also gets a bookmark wrapped around the replaced text)
function Replace-JournalMention($searchStart) {
    $found = $d.Range($searchStart, $searchStart)
    $ok = $found.Find.Execute("the Quarterly Journal of Experimental Psychology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) { return $null }

    $theLen = 4  # length of "the "
    $theStart = $found.Start
    $theEnd = $found.Start + $theLen
    $journalStart = $theEnd
    $journalEnd = $found.End

    # Remove the leading "the " run entirely.
    $theRange = $d.Range($theStart, $theEnd)
    $theRange.Text = ""

    # The journal-name text shifted left by $theLen once "the " was removed.
    $journalRange = $d.Range($journalStart - $theLen, $journalEnd - $theLen)
    $journalRange.Text = "Metacognition and Learning"
    return $journalRange
}

$journal1 = Replace-JournalMention 0
if ($journal1 -ne $null) {
    $d.Bookmarks.Add("_Hlk84672896", $journal1) | Out-Null
    $journal2 = Replace-JournalMention $journal1.End
}
